# Update (Removed Auto Arima)
# Refresh the per-week forecast figures on "Forecast Comparison" and the
# roll-up stats on "Summary" now that the Auto-ARIMA model has been
# dropped from the ensemble (Prophet / Amazon Mean / P70 / P80 / P90).

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Columns: C=Prophet Forecast, D=Amazon Mean Forecast, E=Amazon P70 Forecast,
#          F=Amazon P80 Forecast, G=Amazon P90 Forecast
# Row layout: row, C, D, E, F, G  (rows 2-17 == weeks W01-W16)
$forecastRows = @(
    @(2,  19, 19, 24, 29, 37),
    @(3,  17, 15, 19, 24, 32),
    @(4,  15, 13, 16, 20, 26),
    @(5,  13, 18, 22, 27, 36),
    @(6,  12, 18, 22, 29, 39),
    @(7,  11, 18, 22, 28, 39),
    @(8,  10, 20, 24, 31, 44),
    @(9,  10, 20, 25, 34, 49),
    @(10,  9, 19, 24, 31, 42),
    @(11,  8, 20, 24, 32, 46),
    @(12,  7, 20, 24, 32, 47),
    @(13,  5, 22, 26, 36, 52),
    @(14,  5, 20, 24, 33, 48),
    @(15,  3, 19, 23, 32, 47),
    @(16,  2, 19, 23, 32, 46),
    @(17,  1, 18, 21, 29, 44)
)

foreach ($r in $forecastRows) {
    $row = $r[0]
    $wsForecast.Cells.Item($row, 3).Value = $r[1]   # C - Prophet Forecast
    $wsForecast.Cells.Item($row, 4).Value = $r[2]   # D - Amazon Mean Forecast
    $wsForecast.Cells.Item($row, 5).Value = $r[3]   # E - Amazon P70 Forecast
    $wsForecast.Cells.Item($row, 6).Value = $r[4]   # F - Amazon P80 Forecast
    $wsForecast.Cells.Item($row, 7).Value = $r[5]   # G - Amazon P90 Forecast
}

# Summary sheet B-column values are stored as text (matching the rest of the
# column), so force text entry rather than letting COM coerce them to
# numbers/dates, then drop back to the Normal style so no stray number
# format sticks to the cell.
function Set-SummaryText {
    param($cellRange, [string]$text)
    $cellRange.NumberFormat = "@"
    $cellRange.Value2 = $text
    $cellRange.Style = "Normal"
}

Set-SummaryText $wsSummary.Range("B9")  "147"         # Total Forecast (16 Weeks)
Set-SummaryText $wsSummary.Range("B10") "107"         # Total Forecast (8 Weeks)
Set-SummaryText $wsSummary.Range("B11") "64"          # Total Forecast (4 Weeks)
Set-SummaryText $wsSummary.Range("B12") "19"          # Max Forecast
Set-SummaryText $wsSummary.Range("B14") "1"           # Min Forecast
Set-SummaryText $wsSummary.Range("B15") "2025-03-23"  # Min Forecast Week

Write-Output "Forecast Comparison + Summary refreshed (Auto-ARIMA removed)"
